$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Append one brand-new row (61017 / "qormas ice") right after the old last row (1390) ---
$ws.Range("A1391").Value = 61017
$ws.Range("B1391").Value = "qormas ice"
$ws.Range("B1391").NumberFormat = "@"

# --- 2. Insert 30 blank rows before row 1340 (pushes old 1340-1391 -> 1370-1421) ---
$ws.Rows("1340:1369").Insert()

# --- 3. Fill the 30 newly inserted rows with the new "jixa"/"jala" texture set ---
$newRows = @(
    @{A=20200; B="jixa-rocks"; Style=$true},
    @{A=20201; B="jixa-grass"; Style=$true},
    @{A=20202; B="jixa-dirt-road"; Style=$true},
    @{A=20203; B="jixa-plank-side"; Style=$false},
    @{A=20204; B="jixa-plank-top"; Style=$false},
    @{A=20205; B="jixa-fountain-side"; Style=$false},
    @{A=20206; B="jixa-fountain-side2"; Style=$false},
    @{A=20207; B="jixa - fireplace"; Style=$false},
    @{A=20208; B="jixa-build-a-1"; Style=$false},
    @{A=20209; B="jixa : build 1 - roof (Thatch)"; Style=$false},
    @{A=20210; B="jixa-build-a-2 (door)"; Style=$false},
    @{A=20211; B="jixa-build-a-2"; Style=$false},
    @{A=20212; B="jixa-build-a-3"; Style=$false},
    @{A=20213; B="jixa: dragon transparen"; Style=$false},
    @{A=20214; B="jixa-build-b-1"; Style=$false},
    @{A=20215; B="jixa-build-b-2"; Style=$false},
    @{A=20216; B="jixa-build-b-3 (door)"; Style=$false},
    @{A=20217; B="jixa-build-b-4"; Style=$false},
    @{A=20218; B="jixa-build-b-5"; Style=$false},
    @{A=20219; B="jixa-build-b-5 (roof)"; Style=$false},
    @{A=20220; B="jixa-build-c-1"; Style=$false},
    @{A=20221; B="jixa-build-c-2"; Style=$false},
    @{A=20222; B="jixa-fance-2"; Style=$false},
    @{A=20223; B="jixa - fireplace - top"; Style=$false},
    @{A=20224; B="jixa - build d 1"; Style=$false},
    @{A=20225; B="jixa - build d 2 (roof)"; Style=$false},
    @{A=20226; B="jixa - build d 3 (door)"; Style=$false},
    @{A=20227; B="jala - ruin - pillar"; Style=$false},
    @{A=20228; B="jala - pillar-top"; Style=$false},
    @{A=20229; B="jala - fance"; Style=$false}

)

$startRow = 1340
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $item = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    if (-not $item.Style) {
        $ws.Cells.Item($r, 2).Style = "Normal"
    }
}

# --- 4. Update the sheet-scoped defined name "texture_index" to the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!texture_index") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$B`$1404"
    }
}

# --- 5. Restore the view/selection state recorded in the saved workbook ---
$excel.Goto($ws.Range("A1403"), $true)
$ws.Range("D1371").Select()
